$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 2 (current row 2 becomes row 5, etc.)
$ws.Rows.Item(2).Resize(3).EntireRow.Insert()

# Fill in the new data rows (2-4) with the early measurements
$newData = @(
    @(100, 10, 9.8, 0.4, 0, -84),
    @(250, 9.9, 9.8, 0.7, 1, -80),
    @(400, 9.9, 9.8, 0.7, 1, -75)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = 2 + $i
    $rowVals = $newData[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

# Update sheet view to match new scroll/selection position
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("A38:D38").Select()
